# Financials update: insert a new "most recent period" column before column D
# on the LECO sheet (Income Statement / Balance Sheet / Cash Flow blocks all
# share the same column layout), shifting the existing D:K data right to
# E:L, and populating the new column D with the latest period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank column before D - this shifts D:K -> E:L (and the used
#    range / dimension) automatically, matching the diff's column shift.
$ws.Range("D:D").Insert()

# 2) Populate the brand-new column D with the new period's values.
#    (Row numbers below match the sheet's existing row layout; rows not
#    listed here - section headers, spacer rows, blank subtotal rows - are
#    left as the blank cells the Insert() already produced.)
$newValues = @(
    @{Row=7;   Value=43465},
    @{Row=8;   Value=3028700},
    @{Row=9;   Value=2000200},
    @{Row=10;  Value=1028500},
    @{Row=12;  Value=54200},
    @{Row=13;  Value=0},
    @{Row=14;  Value=25300},
    @{Row=15;  Value=0},
    @{Row=17;  Value=2653100},
    @{Row=18;  Value=375500},
    @{Row=20;  Value=17600},
    @{Row=21;  Value=465500},
    @{Row=22;  Value=24500},
    @{Row=23;  Value=368700},
    @{Row=24;  Value=81300},
    @{Row=25;  Value=0},
    @{Row=26;  Value=287400},
    @{Row=27;  Value=287500},
    @{Row=28;  Value=0},
    @{Row=29;  Value=-400},
    @{Row=30;  Value=0},
    @{Row=31;  Value=0},
    @{Row=32;  Value=-17600},
    @{Row=33;  Value=287100},
    @{Row=34;  Value=0},
    @{Row=35;  Value=287100},
    @{Row=38;  Value=43465},
    @{Row=41;  Value=358800},
    @{Row=42;  Value=0},
    @{Row=43;  Value=396900},
    @{Row=44;  Value=361800},
    @{Row=45;  Value=120200},
    @{Row=46;  Value=1237800},
    @{Row=47;  Value="NA"},
    @{Row=48;  Value=478800},
    @{Row=49;  Value=429200},
    @{Row=50;  Value=0},
    @{Row=51;  Value=0},
    @{Row=52;  Value=204000},
    @{Row=53;  Value=0},
    @{Row=54;  Value=2349800},
    @{Row=57;  Value=268600},
    @{Row=58;  Value=100},
    @{Row=59;  Value=269500},
    @{Row=60;  Value=538200},
    @{Row=61;  Value=702500},
    @{Row=62;  Value=221500},
    @{Row=63;  Value=0},
    @{Row=64;  Value=0},
    @{Row=65;  Value=0},
    @{Row=66;  Value=1462900},
    @{Row=68;  Value=0},
    @{Row=69;  Value=0},
    @{Row=70;  Value=0},
    @{Row=71;  Value=0},
    @{Row=72;  Value=2564400},
    @{Row=73;  Value=0},
    @{Row=74;  Value=0},
    @{Row=75;  Value=0},
    @{Row=76;  Value=886900},
    @{Row=77;  Value=0},
    @{Row=80;  Value=43465},
    @{Row=81;  Value=287100},
    @{Row=83;  Value=72300},
    @{Row=84;  Value=0},
    @{Row=85;  Value=0},
    @{Row=86;  Value=0},
    @{Row=87;  Value=0},
    @{Row=88;  Value=0},
    @{Row=89;  Value=329200},
    @{Row=91;  Value=-71200},
    @{Row=92;  Value=0},
    @{Row=93;  Value=0},
    @{Row=94;  Value=20800},
    @{Row=96;  Value=-102100},
    @{Row=97;  Value=0},
    @{Row=98;  Value=0},
    @{Row=99;  Value=0},
    @{Row=100; Value=-302100},
    @{Row=101; Value=-15700},
    @{Row=102; Value=32100}
)

foreach ($item in $newValues) {
    $ws.Cells.Item($item.Row, 4).Value = $item.Value
}

# 3) The Insert() above left the new column D cells with the default style;
#    copy the (now shifted) column E's formatting onto D so number formats
#    (date header row / plain numbers) match the rest of each row, exactly
#    like using Excel's "Insert Copied Cells" / paste-formatting workflow.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
